$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted before the existing row 322,
# pushing every subsequent record (322-400) down by one row and growing
# the used range to A1:R401.
$ws.Rows.Item(322).Insert()

$ws.Cells.Item(322, 1).Value = 8
$ws.Cells.Item(322, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(322, 3).Value = "Coquimbo"
$ws.Cells.Item(322, 4).Value = 45015
$ws.Cells.Item(322, 5).Value = 4
$ws.Cells.Item(322, 6).Value = 100112021
$ws.Cells.Item(322, 7).Value = "Ají"
$ws.Cells.Item(322, 8).Value = "Inferno"
$ws.Cells.Item(322, 9).Value = "Primera"
$ws.Cells.Item(322, 10).Value = 400
$ws.Cells.Item(322, 11).Value = 11000
$ws.Cells.Item(322, 12).Value = 12000
$ws.Cells.Item(322, 13).Value = 11500
$ws.Cells.Item(322, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(322, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(322, 16).Value = 767
$ws.Cells.Item(322, 17).Value = 15
$ws.Cells.Item(322, 18).Value = "Hortaliza"
